# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly scraped counts (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 761
$wsExhibit.Range("F3").Value = 17
$wsExhibit.Range("F5").Value = 32
$wsExhibit.Range("F6").Value = 260
$wsExhibit.Range("F7").Value = 3562
$wsExhibit.Range("F8").Value = 73
$wsExhibit.Range("F9").Value = 4197
$wsExhibit.Range("F10").Value = 484
$wsExhibit.Range("F11").Value = 1048
$wsExhibit.Range("F12").Value = 50

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 761
$wsAll.Range("F3").Value = 17
$wsAll.Range("F5").Value = 32
$wsAll.Range("F7").Value = 260
$wsAll.Range("F8").Value = 3562
$wsAll.Range("F9").Value = 73
$wsAll.Range("F10").Value = 4198
$wsAll.Range("F11").Value = 484
$wsAll.Range("F12").Value = 1048
$wsAll.Range("F13").Value = 50
